# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" quarter sheet (inserted right after "总计"), fills the
# summary ("总计") sheet with the new quarter's totals as its new first data
# row (existing rows shift down), matching the commit "feat: add 2022-Q3 data".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by copying the "2022-Q2" sheet (so it
#    inherits the same column widths / header & index-column styling), then
#    placing it right after "总计" and renaming it.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy([System.Reflection.Missing]::Value, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2) Fill in the 2022-Q3 fund holdings data (23 funds).
#    Columns: A index(0-based,n) B code(text) C name(text) D size(text)
#             E position(text) F pct(text) G value(text) H rank(n)
# ---------------------------------------------------------------------------
$q3Data = @(
    @("011866", "广发价值增长混合A", "15.73", "94.64", "8.13", "1.2788", 2),
    @("070019", "嘉实价值优势混合A", "22.89", "91.19", "4.86", "1.1125", 10),
    @("002624", "广发优企精选灵活配置混合A", "11.48", "94.33", "8.80", "1.0102", 2),
    @("270025", "广发行业领先混合A", "9.95", "94.64", "8.47", "0.8428", 2),
    @("970016", "中信建投价值增长混合A", "14.00", "86.79", "5.82", "0.8148", 2),
    @("001878", "嘉实沪港深精选股票", "22.02", "91.33", "3.38", "0.7443", 10),
    @("160726", "嘉实瑞享定期开放灵活配置混合", "12.64", "83.43", "4.98", "0.6295", 7),
    @("009138", "嘉实瑞成两年持有期混合A", "10.97", "90.59", "3.56", "0.3905", 8),
    @("011427", "广发价值驱动混合A", "2.64", "94.69", "8.50", "0.2244", 2),
    @("970017", "中信建投价值增长混合C", "2.71", "86.79", "5.82", "0.1577", 2),
    @("009139", "嘉实瑞成两年持有期混合C", "2.99", "90.59", "3.56", "0.1064", 8),
    @("000747", "广发逆向策略灵活配置混合A", "1.00", "94.25", "8.67", "0.0867", 2),
    @("210002", "金鹰红利价值混合A", "1.19", "61.62", "6.23", "0.0741", 1),
    @("011867", "广发价值增长混合C", "0.72", "94.64", "8.13", "0.0585", 2),
    @("011428", "广发价值驱动混合C", "0.42", "94.69", "8.50", "0.0357", 2),
    @("011765", "兴银高端制造混合A", "0.57", "92.99", "3.94", "0.0225", 2),
    @("016563", "金鹰红利价值混合C", "0.34", "61.62", "6.23", "0.0212", 1),
    @("011766", "兴银高端制造混合C", "0.34", "92.99", "3.94", "0.0134", 2),
    @("010021", "广发优企精选灵活配置混合C", "0.07", "94.33", "8.80", "0.0062", 2),
    @("016169", "嘉实价值优势混合C", "0.05", "91.19", "4.86", "0.0024", 10),
    @("005146", "兴银丰润灵活配置混合", "0.04", "92.81", "4.51", "0.0018", 3),
    @("011758", "广发逆向策略灵活配置混合C", "0.02", "94.25", "8.67", "0.0017", 2),
    @("960001", "广发行业领先混合H", "0.01", "94.64", "8.47", "0.0008", 2)
)

$lastRow = 1 + $q3Data.Length   # header row + 23 data rows = row 24

# Add the 3 extra rows this quarter needs (old 2022-Q2 sheet only had 21 rows,
# we need 24); copy formats down from the last existing data row (21) first
# so new rows 22-24 pick up the plain (unstyled) body formatting.
$q3Sheet.Range("A21:H21").Copy()
$q3Sheet.Range("A22:H24").PasteSpecial(-4122)

# Force columns B:G to be stored as text (the source data is textual, e.g.
# "011866", "15.73"), then reset the style back to Normal so only the text
# type sticks and no stray number-format style is left behind.
$bodyRange = $q3Sheet.Range("B2:G" + $lastRow)
$bodyRange.NumberFormat = "@"

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = 2 + $i
    $row = $q3Data[$i]
    $q3Sheet.Cells.Item($r, 1).Value = $i            # A: 0-based index
    $q3Sheet.Cells.Item($r, 2).Value = $row[0]        # B: 基金代码
    $q3Sheet.Cells.Item($r, 3).Value = $row[1]        # C: 基金名称
    $q3Sheet.Cells.Item($r, 4).Value = $row[2]        # D: 基金规模
    $q3Sheet.Cells.Item($r, 5).Value = $row[3]        # E: 股票总仓位
    $q3Sheet.Cells.Item($r, 6).Value = $row[4]        # F: 仓位占比
    $q3Sheet.Cells.Item($r, 7).Value = $row[5]        # G: 持有市值(亿元)
    $q3Sheet.Cells.Item($r, 8).Value = $row[6]        # H: 仓位排名 (number)
}

$bodyRange.Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: insert the new 2022-Q3 totals as the
#    new first data row (row 2), pushing every other quarter down by one row.
# ---------------------------------------------------------------------------
$totalData = @(
    @("2022-Q3", 23, 7.64),
    @("2022-Q2", 20, 6.8),
    @("2022-Q1", 14, 4.33),
    @("2021-Q4", 6, 1.75),
    @("2021-Q3", 6, 2.06),
    @("2021-Q2", 9, 2.74),
    @("2021-Q1", 11, 3.14),
    @("2020-Q4", 11, 4.87)
)

# Extend the formatted index column (A) down into the new row 9, copying the
# format from row 8 so it keeps the bold/centered/bordered look.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = 2 + $i
    $row = $totalData[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i       # A: 0-based index
    $totalSheet.Cells.Item($r, 2).Value = $row[0]  # B: 日期 (quarter label)
    $totalSheet.Cells.Item($r, 3).Value = $row[1]  # C: 持有数量(只)
    $totalSheet.Cells.Item($r, 4).Value = $row[2]  # D: 持有市值(亿元)
}

Write-Host "2022-Q3 sheet added and 总计 updated."
